# The commit swaps the presentation's theme colour palette from the
# "Integral" scheme (ppt/theme/theme1.xml, used by the slide master) to
# the stock Office "Office Theme" palette. The font scheme and the
# format scheme (fills/lines/effects) in theme1.xml already match the
# default Office scheme byte-for-byte, so only the 12 colour-scheme
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) actually need
# to change.
#
# PowerPoint's object model doesn't give VBA/COM a way to rewrite a
# <a:theme>/<a:clrScheme> wholesale, but Slide.ThemeColorScheme exposes
# exactly those twelve slots (in that order) as settable RGBColor
# objects, and writing to them edits the shared theme part used by the
# slide master (ppt/theme/theme1.xml) rather than a per-slide override.

function Convert-HexToRgbValue($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colours, in Slide.ThemeColorScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Length; $i++) {
    $tcs.Item($i).RGB = Convert-HexToRgbValue $officeThemeHex[$i - 1]
}
